$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '58.333.92'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -3.71%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.715.54'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -6.49%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '501.68'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -5.00%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '140.14'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.83%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.528'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -5.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.728.04'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -6.23%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.04'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.06%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.104'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.12%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.347'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.69%  '
$ws.Range('E13').Value = '  +0.97%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.195.61'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -6.09%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '58.474.07'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.39%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.63'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -4.31%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.726.12'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -6.21%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0000135'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -4.91%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.74'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -5.73%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.94'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -5.96%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '341.88'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -6.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.24'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -4.95%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.996'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.65'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '62.81'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.26%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.426'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -5.64%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.172'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -4.98%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.995'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.48'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.65%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0826'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -4.25%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.998'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '19.14'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.59'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.90%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '151.14'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.42'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.94%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.17'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.03%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.944'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -6.19%  '
$ws.Range('E38').Value = '  -6.36%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '35.90'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -5.53%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.39'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -7.53%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.54'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.63%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.184.22'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -6.31%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.997'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0556'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.00%  '
$ws.Range('E45').Value = '  -6.59%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '18.94'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -8.56%  '
$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.37'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.75'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -5.61%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0227'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.36%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0883'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -5.61%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '18.03'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.15%  '
